$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct the amplitude (damping) reading in C2 and propagate it down
#    the column as a "=cell above" reference chain (C3:C11), replacing the
#    old static duplicated values.
$ws.Range("C2").Value = 0.141
for ($r = 3; $r -le 11; $r++) {
    $prev = $r - 1
    $ws.Range("C$r").Formula = "=C$prev"
}

# 2. Correct the fm reading in M2 (fed by the existing M3:M11 = M{r-1} chain).
$ws.Range("M2").Value = 0.09

# 3. Fix the gravitational acceleration constant used throughout the T
#    column (9.8 -> 9.806), row by row.
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("T$r").Formula = "=4000*N$r*9.806/(P$r^2*PI())"
}

# 4. O and Q were static duplicated values; turn them into reference
#    chains mirroring the existing pattern used by L, M, N, P, R columns.
for ($r = 3; $r -le 11; $r++) {
    $prev = $r - 1
    $ws.Range("O$r").Formula = "=O$prev"
    $ws.Range("Q$r").Formula = "=Q$prev"
}

# 5. Update the active selection to K2, matching the saved sheet view.
$ws.Range("K2").Select()

$wb.Application.Calculate()
